$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Spp1-Itgb3 NATMI results (per Dr Hou's advice): sending/target cluster
# pairs recomputed across ECs, FAPs, M2, sCs (4x4 = 16 rows), replacing the
# previous 3x4 (12-row) table.
$data = @(
    @("ECs", "Spp1", "Itgb3", "ECs", 3, 1, 463.732605, 1391.197815, 0.3632113435366598, 0.3632113435366598, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 1150.74751282638, 10356.72761543742, 0.08517556921116826, 0.08517556921116824),
    @("ECs", "Spp1", "Itgb3", "FAPs", 3, 1, 463.732605, 1391.197815, 0.3632113435366598, 0.3632113435366598, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 3359.847990526335, 30238.63191473702, 0.2486878849324593, 0.2486878849324592),
    @("ECs", "Spp1", "Itgb3", "M2", 3, 1, 463.732605, 1391.197815, 0.3632113435366598, 0.3632113435366598, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 92.83354815220498, 835.5019333698449, 0.006871316442245014, 0.006871316442245012),
    @("ECs", "Spp1", "Itgb3", "sCs", 3, 1, 463.732605, 1391.197815, 0.3632113435366598, 0.3632113435366598, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 303.665248844475, 2732.987239600275, 0.02247657295078734, 0.02247657295078733),
    @("FAPs", "Spp1", "Itgb3", "ECs", 3, 1, 12.24662533333333, 36.739876, 0.009591978638444229, 0.009591978638444227, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 30.38987013399644, 273.508831205968, 0.002249385254423894, 0.002249385254423893),
    @("FAPs", "Spp1", "Itgb3", "FAPs", 3, 1, 12.24662533333333, 36.739876, 0.009591978638444229, 0.009591978638444227, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 88.72958052395067, 798.566224715556, 0.006567550607546652, 0.006567550607546649),
    @("FAPs", "Spp1", "Itgb3", "M2", 3, 1, 12.24662533333333, 36.739876, 0.009591978638444229, 0.009591978638444227, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 2.451623350020889, 22.064610150188, 0.0001814632766978886, 0.0001814632766978885),
    @("FAPs", "Spp1", "Itgb3", "sCs", 3, 1, 12.24662533333333, 36.739876, 0.009591978638444229, 0.009591978638444227, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 8.019437255984444, 72.17493530386001, 0.0005935794997757965, 0.0005935794997757962),
    @("M2", "Spp1", "Itgb3", "ECs", 3, 1, 689.7685036666667, 2069.305511, 0.5402504422695089, 0.5402504422695089, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 1711.653184318127, 15404.87865886315, 0.126692460892941, 0.1266924608929409),
    @("M2", "Spp1", "Itgb3", "FAPs", 3, 1, 689.7685036666667, 2069.305511, 0.5402504422695089, 0.5402504422695089, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 4997.529386515333, 44977.76447863799, 0.3699051315787697, 0.3699051315787695),
    @("M2", "Spp1", "Itgb3", "M2", 3, 1, 689.7685036666667, 2069.305511, 0.5402504422695089, 0.5402504422695089, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 138.0831472891881, 1242.748325602693, 0.01022058317548646, 0.01022058317548646),
    @("M2", "Spp1", "Itgb3", "sCs", 3, 1, 689.7685036666667, 2069.305511, 0.5402504422695089, 0.5402504422695089, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 451.6799596418706, 4065.119636776835, 0.03343226662231192, 0.03343226662231191),
    @("sCs", "Spp1", "Itgb3", "ECs", 3, 1, 111.00921, 333.02763, 0.08694623555538696, 0.08694623555538696, 3, 1, 2.481489333333333, 7.444467999999999, 0.2345069082418988, 0.2345069082418987, 275.4681705167599, 2479.213534650839, 0.02038949288336564, 0.02038949288336564),
    @("sCs", "Spp1", "Itgb3", "FAPs", 3, 1, 111.00921, 333.02763, 0.08694623555538696, 0.08694623555538696, 3, 1, 7.245227, 21.735681, 0.6846919551326144, 0.6846919551326142, 804.2869255406699, 7238.58232986603, 0.05953138801383873, 0.05953138801383871),
    @("sCs", "Spp1", "Itgb3", "M2", 3, 1, 111.00921, 333.02763, 0.08694623555538696, 0.08694623555538696, 2, 0.6666666666666666, 0.2001876666666667, 0.600563, 0.01891823194544989, 0.01891823194544989, 22.22267472841, 200.00407255569, 0.001644869051020533, 0.001644869051020532),
    @("sCs", "Spp1", "Itgb3", "sCs", 3, 1, 111.00921, 333.02763, 0.08694623555538696, 0.08694623555538696, 3, 1, 0.6548283333333333, 1.964485, 0.06188290468003712, 0.06188290468003711, 72.69197596895, 654.22778372055, 0.005380485607162065, 0.005380485607162064)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
